# Auto-generated edit script applying scheduled runner updates to Sheets workbook
# Updates currentAveragePrice / Leve profit calculations across ALC, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

# ALC row 11
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 499
$ws.Range("I11").Value = 499
$ws.Range("K11").Value = 499
$ws.Range("M11").Value = -359

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4624896
$ws.Range("J17").Value = 4624896
$ws.Range("L17").Value = 13874688
$ws.Range("N17").Value = -13875024

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 824.7273
$ws.Range("I28").Value = 376.05264
$ws.Range("K28").Value = 376.05264
$ws.Range("M28").Value = 108.94736

# ALC row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2161.8333
$ws.Range("I38").Value = 269.6
$ws.Range("J38").Value = 3513.4285
$ws.Range("K38").Value = 808.8000000000001
$ws.Range("L38").Value = 10540.2855
$ws.Range("M38").Value = -436.8000000000001
$ws.Range("N38").Value = -11284.2855

# ALC row 42
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 2917.611
$ws.Range("I42").Value = 1591.375
$ws.Range("J42").Value = 3978.6
$ws.Range("K42").Value = 4774.125
$ws.Range("L42").Value = 11935.8
$ws.Range("M42").Value = -4544.125
$ws.Range("N42").Value = -12395.8

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6386.375
$ws.Range("J43").Value = 6033.6665
$ws.Range("L43").Value = 6033.6665
$ws.Range("N43").Value = -6171.6665

# ALC row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 6817.5
$ws.Range("I106").Value = 7381.2
$ws.Range("J106").Value = 3999
$ws.Range("K106").Value = 7381.2
$ws.Range("L106").Value = 3999
$ws.Range("M106").Value = -6750.2
$ws.Range("N106").Value = -5261

# ALC row 113
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 12929.2
$ws.Range("I113").Value = 17067
$ws.Range("K113").Value = 17067
$ws.Range("M113").Value = -13813

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 6986.25
$ws.Range("I116").Value = 6841.4287
$ws.Range("J116").Value = 8000
$ws.Range("K116").Value = 6841.4287
$ws.Range("L116").Value = 8000
$ws.Range("M116").Value = -3399.4287
$ws.Range("N116").Value = -14884

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 17061.785
$ws.Range("J137").Value = 28937
$ws.Range("L137").Value = 86811
$ws.Range("N137").Value = -91911

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2428.8555
$ws.Range("I138").Value = 1954.8235
$ws.Range("J138").Value = 2550.9546
$ws.Range("K138").Value = 5864.470499999999
$ws.Range("L138").Value = 7652.8638
$ws.Range("M138").Value = -724.4704999999994
$ws.Range("N138").Value = -17932.8638

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6962484
$ws.Range("I20").Value = 13894070
$ws.Range("K20").Value = 13894070
$ws.Range("M20").Value = -13893823

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 47621884
$ws.Range("I86").Value = 3570.5454
$ws.Range("J86").Value = 100002024
$ws.Range("K86").Value = 3570.5454
$ws.Range("L86").Value = 100002024
$ws.Range("M86").Value = -2447.5454
$ws.Range("N86").Value = -100004270

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 47621884
$ws.Range("I89").Value = 3570.5454
$ws.Range("J89").Value = 100002024
$ws.Range("K89").Value = 17852.727
$ws.Range("L89").Value = 500010120
$ws.Range("M89").Value = -12236.727
$ws.Range("N89").Value = -500021352

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 19094.94
$ws.Range("I99").Value = 18155.428
$ws.Range("K99").Value = 18155.428
$ws.Range("M99").Value = -16657.428

# BSM row 109
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 59994.5

# CRP row 53
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value = 55000
$ws.Range("J53").Value = 55000
$ws.Range("L53").Value = 55000
$ws.Range("N53").Value = -56214

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 11141.88
$ws.Range("I58").Value = 4070.0857
$ws.Range("J58").Value = 27642.732
$ws.Range("K58").Value = 4070.0857
$ws.Range("L58").Value = 27642.732
$ws.Range("M58").Value = -3867.0857
$ws.Range("N58").Value = -28048.732

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2324.611
$ws.Range("I132").Value = 1990.7646
$ws.Range("J132").Value = 8000
$ws.Range("K132").Value = 5972.293799999999
$ws.Range("L132").Value = 24000
$ws.Range("M132").Value = -3442.293799999999
$ws.Range("N132").Value = -29060

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 11141.88
$ws.Range("I136").Value = 4070.0857
$ws.Range("J136").Value = 27642.732
$ws.Range("K136").Value = 12210.2571
$ws.Range("L136").Value = 82928.196
$ws.Range("M136").Value = -9660.257100000001
$ws.Range("N136").Value = -88028.196

# CUL row 17
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 500.7143
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

# CUL row 22
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 6154.9165
$ws.Range("J22").Value = 6982
$ws.Range("L22").Value = 20946
$ws.Range("N22").Value = -21284

# CUL row 27
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 6154.9165
$ws.Range("J27").Value = 6982
$ws.Range("L27").Value = 20946
$ws.Range("N27").Value = -21150

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1131.5294
$ws.Range("I113").Value = 1102.5
$ws.Range("J113").Value = 1157.3334
$ws.Range("K113").Value = 3307.5
$ws.Range("L113").Value = 3472.0002
$ws.Range("M113").Value = -1137.5
$ws.Range("N113").Value = -7812.0002

# CUL row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 10102012
$ws.Range("I129").Value = 1170.5
$ws.Range("J129").Value = 30303696
$ws.Range("K129").Value = 3511.5
$ws.Range("L129").Value = 90911088
$ws.Range("M129").Value = 1488.5
$ws.Range("N129").Value = -90921088

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1474.45
$ws.Range("J131").Value = 1478.9395
$ws.Range("L131").Value = 4436.818499999999
$ws.Range("N131").Value = -14516.8185

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 16740.834
$ws.Range("I80").Value = 16699
$ws.Range("J80").Value = 16799.4
$ws.Range("K80").Value = 16699
$ws.Range("L80").Value = 16799.4
$ws.Range("M80").Value = -15701
$ws.Range("N80").Value = -18795.4

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 16740.834
$ws.Range("I83").Value = 16699
$ws.Range("J83").Value = 16799.4
$ws.Range("K83").Value = 83495
$ws.Range("L83").Value = 83997
$ws.Range("M83").Value = -78503
$ws.Range("N83").Value = -93981

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 66668576
$ws.Range("I16").Value = 83335510
$ws.Range("K16").Value = 83335510
$ws.Range("M16").Value = -83335340

# LTW row 32
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 6825
$ws.Range("I32").Value = 2433.3333
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 2433.3333
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -2116.3333
$ws.Range("N32").Value = -20634

# LTW row 34
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("N34").Value = 0
$ws.Range("L34").ClearContents()
$ws.Range("M34").ClearContents()

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3735.4546

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 3043.1
$ws.Range("I82").Value = 3290.1538
$ws.Range("J82").Value = 2584.2856
$ws.Range("K82").Value = 3290.1538
$ws.Range("L82").Value = 2584.2856
$ws.Range("M82").Value = -2929.1538
$ws.Range("N82").Value = -3306.2856

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 3043.1
$ws.Range("I85").Value = 3290.1538
$ws.Range("J85").Value = 2584.2856
$ws.Range("K85").Value = 3290.1538
$ws.Range("L85").Value = 2584.2856
$ws.Range("M85").Value = -2042.1538
$ws.Range("N85").Value = -5080.2856

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 5385.364
$ws.Range("I100").Value = 5939.222
$ws.Range("K100").Value = 5939.222
$ws.Range("M100").Value = -5398.222

# LTW row 109
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H109").Value = 12394
$ws.Range("J109").Value = 12394
$ws.Range("L109").Value = 12394
$ws.Range("N109").Value = -15168

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2119333.8
$ws.Range("J132").Value = 4471768
$ws.Range("L132").Value = 13415304
$ws.Range("N132").Value = -13420364

# WVR row 42
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 24750
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()

# WVR row 43
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 19700
$ws.Range("I43").Value = 12833.333
$ws.Range("J43").Value = 30000
$ws.Range("K43").Value = 12833.333
$ws.Range("L43").Value = 30000
$ws.Range("M43").Value = -12684.333
$ws.Range("N43").Value = -30298

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 26538.154
$ws.Range("I62").Value = 23249.625
$ws.Range("J62").Value = 31799.8
$ws.Range("K62").Value = 23249.625
$ws.Range("L62").Value = 31799.8
$ws.Range("M62").Value = -22625.625
$ws.Range("N62").Value = -33047.8

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 26538.154
$ws.Range("I65").Value = 23249.625
$ws.Range("J65").Value = 31799.8
$ws.Range("K65").Value = 116248.125
$ws.Range("L65").Value = 158999
$ws.Range("M65").Value = -113128.125
$ws.Range("N65").Value = -165239

# WVR row 100
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 600
$ws.Range("J100").Value = 636.7273
$ws.Range("L100").Value = 1273.4546
$ws.Range("N100").Value = -2355.4546

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2529
$ws.Range("I107").Value = 2529
$ws.Range("K107").Value = 7587
$ws.Range("M107").Value = -5667

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 13069.26
$ws.Range("I132").Value = 5089.026
$ws.Range("K132").Value = 15267.078
$ws.Range("M132").Value = -12737.078

